$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Ste"
$ws.Range("B1").Value = "Ets"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Four"
$ws.Range("E1").Value = "Lieu"
$ws.Range("F1").Value = "Depot"
$ws.Range("G1").Value = "Article"
$ws.Range("H1").Value = "Qte"
$ws.Range("I1").Value = "Un"
$ws.Range("J1").Value = "ORI"

$headerRange = $ws.Range("A1:J1")
$headerRange.Interior.Color = 65535
$headerRange.HorizontalAlignment = -4108

$null = $ws.Range("G20").Select()
